$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'329.31"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'1.52%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'41.36"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'4.53%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.618"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.28%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08198"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'2.13%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'8.753"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'1.55%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'2.005"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-0.80%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'4.498"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'0.18%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.987"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'1.94%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.9232"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'0.04%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1279"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'3.03%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.1952"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-0.95%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09276"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'1.60%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.03848"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'7.63%"
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'0.86%"
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'0.66%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.006141"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.51%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D19").Value = "'3.446"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'2.93%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.3479"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'0.01%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'8.223"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-5.83%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.1365"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'0.91%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.2659"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'10.18%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.04396"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.33%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'-0.28%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004314"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-6.33%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0001201"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'-2.56%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.02759"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'10.60%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.05457"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'1.97%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007801"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'4.61%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1419"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'1.16%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.008936"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-9.75%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002171"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'2.46%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.01145"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'2.10%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006772"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'0.23%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'-0.20%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.003191"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'7.10%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.002278"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-0.09%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002099"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.20%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0001999"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.20%"
$ws.Range("E51").Style = "Normal"
